# ============================================================================
# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu
#
# 1. "Đơn sale phụ": add two new order rows (629, 631), recompute the
#    "Tổng" (total) row which shifts down from row 3 to row 5.
# 2. Insert a brand-new sheet "Đơn thu nợ" (debt-collection orders) between
#    "Đơn sale phụ" and "Lương".
# 3. "Lương": add a "Chiết khấu thu nợ" line per cơ sở + refresh all the
#    downstream totals that change because of the new sale/thu nợ rows.
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-Text {
    param($range, $value)
    # Force text interpretation so date-looking strings ("08-05-2024") are
    # not silently reinterpreted as date serials.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ----------------------------------------------------------------------
# 1) "Đơn sale phụ" — insert rows 3 & 4, move + refresh the Tổng row (->5)
# ----------------------------------------------------------------------
$wsPhu = $wb.Worksheets.Item("Đơn sale phụ")

$newOrders = @(
    @{ Row = 3; A = "HD-LUXURY"; B = 629; C = "08-05-2024"; D = "CẦN THƠ"; E = "Phạm Thị Trúc Lài"; F = "CTV"; G = "Phun mày"; H = 500000; I = "Đỗ Thị Huyền Trân"; J = 500000; K = 1000000; L = 1000000; M = 0.02; N = 10000 },
    @{ Row = 4; A = "HD-LUXURY"; B = 631; C = "08-06-2024"; D = "CẦN THƠ"; E = "Võ Thị Thuỳ Trang";  F = "CTV"; G = "Phun mày"; H = 500000; I = "Đỗ Thị Huyền Trân"; J = 1000000; K = 1500000; L = 1500000; M = 0.02; N = 20000 }
)

foreach ($o in $newOrders) {
    $r = $o.Row
    $wsPhu.Cells.Item($r, 1).Value = $o.A
    $wsPhu.Cells.Item($r, 2).Value = $o.B
    Set-Text $wsPhu.Cells.Item($r, 3) $o.C
    $wsPhu.Cells.Item($r, 4).Value = $o.D
    $wsPhu.Cells.Item($r, 5).Value = $o.E
    $wsPhu.Cells.Item($r, 6).Value = $o.F
    $wsPhu.Cells.Item($r, 7).Value = $o.G
    $wsPhu.Cells.Item($r, 8).Value = $o.H
    $wsPhu.Cells.Item($r, 9).Value = $o.I
    $wsPhu.Cells.Item($r, 10).Value = $o.J
    $wsPhu.Cells.Item($r, 11).Value = $o.K
    $wsPhu.Cells.Item($r, 12).Value = $o.L
    $wsPhu.Cells.Item($r, 13).Value = $o.M
    $wsPhu.Cells.Item($r, 14).Value = $o.N
}

# "Tổng" row moves from row 3 down to row 5, with refreshed aggregates.
$wsPhu.Cells.Item(5, 1).Value = "Tổng"
$wsPhu.Cells.Item(5, 2).Value = 3
$wsPhu.Cells.Item(5, 3).Value = ""
$wsPhu.Cells.Item(5, 4).Value = ""
$wsPhu.Cells.Item(5, 5).Value = ""
$wsPhu.Cells.Item(5, 6).Value = ""
$wsPhu.Cells.Item(5, 7).Value = ""
$wsPhu.Cells.Item(5, 8).Value = 1000000
$wsPhu.Cells.Item(5, 9).Value = ""
$wsPhu.Cells.Item(5, 10).Value = 7500000
$wsPhu.Cells.Item(5, 11).Value = 8500000
$wsPhu.Cells.Item(5, 12).Value = 8500000
$wsPhu.Cells.Item(5, 13).Value = 0
$wsPhu.Cells.Item(5, 14).Value = 270000

Write-Host "Updated 'Đơn sale phụ'"

# ----------------------------------------------------------------------
# 2) Insert a brand-new "Đơn thu nợ" sheet, right after "Đơn sale phụ"
# ----------------------------------------------------------------------
$wsNo = $wb.Worksheets.Add($null, $wsPhu)
$wsNo.Name = "Đơn thu nợ"

$noHeaders = @(
    "Tiền tố", "Mã đơn thu nợ", "Lượng thu", "Ngày thu", "Cơ sở", "Đơn nợ",
    "Tên dịch vụ", "Khách hàng", "Nguồn khách", "Sale chính", "Đơn giá gốc",
    "Sale phụ", "Upsale", "Đơn giá", "Đã thanh toán", "Bác sĩ 1", "Bác sĩ 2",
    "Tỉ lệ chiết khấu sale chính", "Chiết khấu sale chính",
    "Tỉ lệ chiết khấu sale phụ", "Chiết khấu sale phụ",
    "Tỉ lệ chiết khấu bác sĩ 1", "Chiết khấu bác sĩ 1",
    "Tỉ lệ chiết khấu bác sĩ 2", "Chiết khấu bác sĩ 2"
)

$c = 1
foreach ($h in $noHeaders) {
    $wsNo.Cells.Item(1, $c).Value = $h
    $c = $c + 1
}

# Row 2 — single debt-collection order
$wsNo.Cells.Item(2, 1).Value = "TN"
$wsNo.Cells.Item(2, 2).Value = 181
$wsNo.Cells.Item(2, 3).Value = 1500000
Set-Text $wsNo.Cells.Item(2, 4) "08-09-2024"
$wsNo.Cells.Item(2, 5).Value = "CẦN THƠ"
$wsNo.Cells.Item(2, 6).Value = "HD-LUXURY-538"
$wsNo.Cells.Item(2, 7).Value = "Nâng mũi"
$wsNo.Cells.Item(2, 8).Value = "Ngô Xuân Nhi"
$wsNo.Cells.Item(2, 9).Value = "Cá nhân"
$wsNo.Cells.Item(2, 10).Value = "Lâm Hoàng Phú"
$wsNo.Cells.Item(2, 11).Value = 10000000
$wsNo.Cells.Item(2, 12).Value = "Đỗ Thị Huyền Trân"
$wsNo.Cells.Item(2, 13).Value = 8000000
$wsNo.Cells.Item(2, 14).Value = 18000000
$wsNo.Cells.Item(2, 15).Value = 11000000
$wsNo.Cells.Item(2, 16).Value = "Lâm Thị Mỹ Hằng"
# Q2 (Bác sĩ 2) intentionally left blank — no second doctor on this order.
$wsNo.Cells.Item(2, 18).Value = 0
$wsNo.Cells.Item(2, 19).Value = 0
$wsNo.Cells.Item(2, 20).Value = 0.04
$wsNo.Cells.Item(2, 21).Value = 40000
$wsNo.Cells.Item(2, 22).Value = 0
$wsNo.Cells.Item(2, 23).Value = 0
$wsNo.Cells.Item(2, 24).Value = 0
$wsNo.Cells.Item(2, 25).Value = 0

# Row 3 — "Tổng" aggregate row
$wsNo.Cells.Item(3, 1).Value = "Tổng"
$wsNo.Cells.Item(3, 2).Value = 1
$wsNo.Cells.Item(3, 3).Value = 1500000
$wsNo.Cells.Item(3, 4).Value = ""
$wsNo.Cells.Item(3, 5).Value = ""
$wsNo.Cells.Item(3, 6).Value = ""
$wsNo.Cells.Item(3, 7).Value = ""
$wsNo.Cells.Item(3, 8).Value = ""
$wsNo.Cells.Item(3, 9).Value = ""
$wsNo.Cells.Item(3, 10).Value = ""
$wsNo.Cells.Item(3, 11).Value = 10000000
$wsNo.Cells.Item(3, 12).Value = ""
$wsNo.Cells.Item(3, 13).Value = 8000000
$wsNo.Cells.Item(3, 14).Value = 18000000
$wsNo.Cells.Item(3, 15).Value = 11000000
$wsNo.Cells.Item(3, 16).Value = ""
$wsNo.Cells.Item(3, 17).Value = ""
$wsNo.Cells.Item(3, 18).Value = 0
$wsNo.Cells.Item(3, 19).Value = 0
$wsNo.Cells.Item(3, 20).Value = 0
$wsNo.Cells.Item(3, 21).Value = 40000
$wsNo.Cells.Item(3, 22).Value = 0
$wsNo.Cells.Item(3, 23).Value = 0
$wsNo.Cells.Item(3, 24).Value = 0
$wsNo.Cells.Item(3, 25).Value = 0

Write-Host "Added 'Đơn thu nợ'"

# ----------------------------------------------------------------------
# 3) "Lương" — insert a "Chiết khấu thu nợ tại <cơ sở>" line for each of
#    the 3 cơ sở, and refresh every total that moved because of it.
# ----------------------------------------------------------------------
$wsLuong = $wb.Worksheets.Item("Lương")

$luongRows = @(
    @{ A = "Danh mục lương"; B = 9 },
    @{ A = "Tổng công tại CẦN THƠ"; B = 6.5 },
    @{ A = "Phụ cấp tại CẦN THƠ"; B = 227500 },
    @{ A = "Lương cơ bản tại CẦN THƠ"; B = 1160714.285714286 },
    @{ A = "Chiết khấu sale chính tại CẦN THƠ"; B = 210000 },
    @{ A = "Chiết khấu sale phụ tại CẦN THƠ"; B = 270000 },
    @{ A = "Đơn 1 bác sĩ tại CẦN THƠ"; B = 0 },
    @{ A = "Đơn 2 bác sĩ tại CẦN THƠ"; B = 0 },
    @{ A = "Công phụ phẫu 1 tại CẦN THƠ"; B = 0 },
    @{ A = "Công phụ phẫu 2 tại CẦN THƠ"; B = 0 },
    @{ A = "Chiết khấu thu nợ tại CẦN THƠ"; B = 40000 },
    @{ A = "Ứng lương tại CẦN THƠ"; B = 0 },
    @{ A = "Tổng công tại LONG XUYÊN"; B = 0 },
    @{ A = "Lương công tác tại LONG XUYÊN"; B = 0 },
    @{ A = "Lương cơ bản tại LONG XUYÊN"; B = $null },
    @{ A = "Chiết khấu sale chính tại LONG XUYÊN"; B = 0 },
    @{ A = "Chiết khấu sale phụ tại LONG XUYÊN"; B = 0 },
    @{ A = "Đơn 1 bác sĩ tại LONG XUYÊN"; B = 0 },
    @{ A = "Đơn 2 bác sĩ tại LONG XUYÊN"; B = 0 },
    @{ A = "Công phụ phẫu 1 tại LONG XUYÊN"; B = 0 },
    @{ A = "Công phụ phẫu 2 tại LONG XUYÊN"; B = 0 },
    @{ A = "Chiết khấu thu nợ tại LONG XUYÊN"; B = 0 },
    @{ A = "Ứng lương tại LONG XUYÊN"; B = 0 },
    @{ A = "Tổng công tại SÓC TRĂNG"; B = 0 },
    @{ A = "Lương công tác tại SÓC TRĂNG"; B = 0 },
    @{ A = "Lương cơ bản tại SÓC TRĂNG"; B = $null },
    @{ A = "Chiết khấu sale chính tại SÓC TRĂNG"; B = 0 },
    @{ A = "Chiết khấu sale phụ tại SÓC TRĂNG"; B = 0 },
    @{ A = "Đơn 1 bác sĩ tại SÓC TRĂNG"; B = 0 },
    @{ A = "Đơn 2 bác sĩ tại SÓC TRĂNG"; B = 0 },
    @{ A = "Công phụ phẫu 1 tại SÓC TRĂNG"; B = 0 },
    @{ A = "Công phụ phẫu 2 tại SÓC TRĂNG"; B = 0 },
    @{ A = "Chiết khấu thu nợ tại SÓC TRĂNG"; B = 0 },
    @{ A = "Ứng lương tại SÓC TRĂNG"; B = 0 },
    @{ A = "Tổng lương tại CẦN THƠ"; B = 1908214.285714286 },
    @{ A = "Tổng lương tại LONG XUYÊN"; B = 0 },
    @{ A = "Tổng lương tại SÓC TRĂNG"; B = 0 },
    @{ A = "Tổng lương tại HỆ THỐNG"; B = 1908214.285714286 }
)

$r = 1
foreach ($row in $luongRows) {
    $wsLuong.Cells.Item($r, 1).Value = $row.A
    if ($null -ne $row.B) {
        $wsLuong.Cells.Item($r, 2).Value = $row.B
    } else {
        # "Lương cơ bản tại LONG XUYÊN/SÓC TRĂNG" stay genuinely blank
        # (not yet computed) — clear any stale value left by the old layout.
        $wsLuong.Cells.Item($r, 2).Clear()
    }
    $r = $r + 1
}

Write-Host "Updated 'Lương'"
